$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 (the "H 72" record), shifting all rows below it up by one.
$ws.Rows.Item(2).Delete()
